$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 2 ("The problem") - shape id=3 "Title 1" - the "We empathized with:"
# bullet box.
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$shEmp = $s2.Shapes.Item(3)
$trEmp = $shEmp.TextFrame.TextRange

# Paragraph 1: "We empathized with:" -> "We " / "empathised" / " with:"
$para1 = $trEmp.Paragraphs(1, 1)
$oldWord = "empathized"
$wordStart = $para1.Start + 3
$runMid = $trEmp.Characters($wordStart, $oldWord.Length)
$runMid.Text = "empathised"
$runMid.Font.Name = "Segoe UI"

# Paragraph 3 ("The fact that there are so many medication options and
# pricing") + paragraph 4 ("Someone who's pregnant, on other meds, needs to
# drive (or all three!)") merge into a single paragraph with new wording,
# split across three runs ("... medication " / "ie" / " pregnancy, ...").
$para3 = $trEmp.Paragraphs(3, 1)
$mergedText = "Someone who has specific considerations when taking medication ie pregnancy, current medications, requires non-drowsy"
$para3Chars = $trEmp.Characters($para3.Start, $para3.Length - 1)
$para3Chars.Text = $mergedText

# Remove the now-redundant old 4th paragraph (its text has been folded into
# paragraph 3 above).
$para4 = $trEmp.Paragraphs(4, 1)
$para4.Delete()

# Split paragraph 3 into its three runs.
$para3b = $trEmp.Paragraphs(3, 1)
$prefix = "Someone who has specific considerations when taking medication "
$ieWord = "ie"
$runIe = $trEmp.Characters($para3b.Start + $prefix.Length, $ieWord.Length)
$runIe.Text = $ieWord
$runIe.Font.Name = "Segoe UI"

# ---------------------------------------------------------------------------
# Slide 2 - shape id=15 "Title 1" - "Our problem statement:" bullet box.
# Last bullet: "In another location e.g. on holiday abroad" ->
# "Numerous medication options and pricing"
# ---------------------------------------------------------------------------
$shProblem = $s2.Shapes.Item(11)
$trProblem = $shProblem.TextFrame.TextRange
$paraLoc = $trProblem.Paragraphs(5, 1)
$paraLocChars = $trProblem.Characters($paraLoc.Start, $paraLoc.Length - 1)
$paraLocChars.Text = "Numerous medication options and pricing"

# ---------------------------------------------------------------------------
# Slide 3 ("The solution") - shape id=22 "Title 1" - "We created" box.
# ---------------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$shCreated = $s3.Shapes.Item(8)

# Resize / reposition the box.
$shCreated.Left = 163.4893798828125
$shCreated.Top = 321.67498779296875
$shCreated.Width = 716.9000854492188
$shCreated.Height = 178.66323852539062

$trCreated = $shCreated.TextFrame.TextRange

# Paragraph 2: "<Insert the solution to the problem>" becomes
# "An AI interface which the user can 'talk' to, stating their symptoms. It
# will ask the user a serious of questions to determine their personal " +
# "circumstances" + " to help determine the most suitable medication."
$para2 = $trCreated.Paragraphs(2, 1)
$para2Chars = $trCreated.Characters($para2.Start, $para2.Length - 1)
$para2Chars.Text = "circumstances"

$para2b = $trCreated.Paragraphs(2, 1)
$prefixSolution = "An AI interface which the user can 'talk' to, stating their symptoms. It will ask the user a serious of questions to determine their personal "
$para2b.InsertBefore($prefixSolution)

$para2c = $trCreated.Paragraphs(2, 1)
$suffixSolution = " to help determine the most suitable medication."
$para2c.InsertAfter($suffixSolution)
